$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the style of the existing header row (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data cells
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 15

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
